$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 629360.1
$ws.Range("I64").Value = 772250.94
$ws.Range("J64").Value = 10166.667
$ws.Range("K64").Value = 772250.94
$ws.Range("L64").Value = 10166.667
$ws.Range("M64").Value = -772002.94
$ws.Range("N64").Value = -10662.667

$ws.Range("H67").Value = 629360.1
$ws.Range("I67").Value = 772250.94
$ws.Range("J67").Value = 10166.667
$ws.Range("K67").Value = 772250.94
$ws.Range("L67").Value = 10166.667
$ws.Range("M67").Value = -771392.94
$ws.Range("N67").Value = -11882.667

$ws.Range("H74").Value = 3280.3
$ws.Range("I74").Value = 3100.5
$ws.Range("J74").Value = 3550
$ws.Range("K74").Value = 3100.5
$ws.Range("L74").Value = 3550
$ws.Range("M74").Value = -2164.5
$ws.Range("N74").Value = -5422

$ws.Range("H77").Value = 3280.3
$ws.Range("I77").Value = 3100.5
$ws.Range("J77").Value = 3550
$ws.Range("K77").Value = 15502.5
$ws.Range("L77").Value = 17750
$ws.Range("M77").Value = -10822.5
$ws.Range("N77").Value = -27110

$ws.Range("H103").Value = 1168.3334
$ws.Range("I103").Value = 1266.6666
$ws.Range("J103").Value = 1070
$ws.Range("K103").Value = 3799.9998
$ws.Range("L103").Value = 3210
$ws.Range("M103").Value = -3213.9998
$ws.Range("N103").Value = -4382

$ws.Range("H111").Value = 1852.7142
$ws.Range("I111").Value = 997.25
$ws.Range("J111").Value = 2993.3333
$ws.Range("K111").Value = 2991.75
$ws.Range("L111").Value = 8979.999899999999
$ws.Range("M111").Value = 75.25
$ws.Range("N111").Value = -15113.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20395.457
$ws.Range("I32").Value = 3296.182
$ws.Range("K32").Value = 3296.182
$ws.Range("M32").Value = -3009.182

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = ""

$ws.Range("H63").Value = 13562.363
$ws.Range("I63").Value = 19411.428
$ws.Range("J63").Value = 3326.5
$ws.Range("K63").Value = 19411.428
$ws.Range("L63").Value = 3326.5
$ws.Range("M63").Value = -18725.428
$ws.Range("N63").Value = -4698.5

$ws.Range("H66").Value = 13562.363
$ws.Range("I66").Value = 19411.428
$ws.Range("J66").Value = 3326.5
$ws.Range("K66").Value = 97057.14
$ws.Range("L66").Value = 16632.5
$ws.Range("M66").Value = -93625.14
$ws.Range("N66").Value = -23496.5

$ws.Range("H88").Value = 5866.3335
$ws.Range("J88").Value = 9224.25
$ws.Range("L88").Value = 9224.25
$ws.Range("N88").Value = -10036.25

$ws.Range("H91").Value = 5866.3335
$ws.Range("J91").Value = 9224.25
$ws.Range("L91").Value = 9224.25
$ws.Range("N91").Value = -12032.25

$ws.Range("H102").Value = 3121.4
$ws.Range("I102").Value = 3293.923
$ws.Range("K102").Value = 3293.923
$ws.Range("M102").Value = -1671.923

$ws.Range("H125").Value = 33000
$ws.Range("J125").Value = 33000
$ws.Range("L125").Value = 33000
$ws.Range("N125").Value = -42840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 24101.4
$ws.Range("I86").Value = 5125
$ws.Range("J86").Value = 100007
$ws.Range("K86").Value = 5125
$ws.Range("L86").Value = 100007
$ws.Range("M86").Value = -4002
$ws.Range("N86").Value = -102253

$ws.Range("H89").Value = 24101.4
$ws.Range("I89").Value = 5125
$ws.Range("J89").Value = 100007
$ws.Range("K89").Value = 25625
$ws.Range("L89").Value = 500035
$ws.Range("M89").Value = -20009
$ws.Range("N89").Value = -511267

$ws.Range("H99").Value = 1736.75
$ws.Range("I99").Value = 1725.8667
$ws.Range("J99").Value = 1900
$ws.Range("K99").Value = 1725.8667
$ws.Range("L99").Value = 1900
$ws.Range("M99").Value = -227.8667
$ws.Range("N99").Value = -4896

$ws.Range("H107").Value = 1246.6666
$ws.Range("I107").Value = 1246.6666
$ws.Range("K107").Value = 1246.6666
$ws.Range("M107").Value = 673.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 15646170
$ws.Range("J99").Value = 2895
$ws.Range("L99").Value = 2895
$ws.Range("N99").Value = -5891

$ws.Range("H126").Value = 15646170
$ws.Range("J126").Value = 2895
$ws.Range("L126").Value = 8685
$ws.Range("N126").Value = -13625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 14706455
$ws.Range("I113").Value = 699
$ws.Range("K113").Value = 2097
$ws.Range("M113").Value = 73

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").Value = ""

$ws.Range("H113").Value = 2104.3635
$ws.Range("I113").Value = 1642.8182
$ws.Range("J113").Value = 2565.9092
$ws.Range("K113").Value = 1642.8182
$ws.Range("L113").Value = 2565.9092
$ws.Range("M113").Value = 527.1818000000001
$ws.Range("N113").Value = -6905.9092

$ws.Range("H126").Value = 2199.2307
$ws.Range("I126").Value = 1822.5
$ws.Range("J126").Value = 2366.6667
$ws.Range("K126").Value = 5467.5
$ws.Range("L126").Value = 7100.000100000001
$ws.Range("M126").Value = -2997.5
$ws.Range("N126").Value = -12040.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").Value = ""

$ws.Range("H61").Value = 12211
$ws.Range("I61").Value = 16650
$ws.Range("J61").Value = 3333
$ws.Range("K61").Value = 16650
$ws.Range("L61").Value = 3333
$ws.Range("M61").Value = -16448
$ws.Range("N61").Value = -3737

$ws.Range("H113").Value = 12211
$ws.Range("I113").Value = 16650
$ws.Range("J113").Value = 3333
$ws.Range("K113").Value = 16650
$ws.Range("L113").Value = 3333
$ws.Range("M113").Value = -14480
$ws.Range("N113").Value = -7673

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = ""

$ws.Range("H117").Value = 32000.5
$ws.Range("J117").Value = 32000.5
$ws.Range("L117").Value = 32000.5
$ws.Range("N117").Value = -41178.5
